$wb = $excel.ActiveWorkbook

# Add a new worksheet (Sheet3) at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Fill column A with the data values
$values = @(10, 12, 2, 2, 1, 2, 3, 2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws3.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Select A8 on the new sheet and make it the active sheet/tab
$ws3.Range("A8").Select()
$ws3.Activate()
